$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks like a plain number (single dot)
# must be forced to Text format first, otherwise Excel auto-converts them
# to a numeric value (losing the exact decimal-string formatting), then
# the number format is reset back to General so no stray style sticks to
# the cell (matches original workbook formatting).

$ws.Range("D2").Value = "46.098.78"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "2.450.87"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "

$ws.Range("D15").Value = "2.839.25"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").Value = "2.465.28"
$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "46.038.90"
$ws.Range("E18").Value = "  +3.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  +2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.95%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "

$ws.Range("E32").Value = "  +5.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.95%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0763"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "127.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("E41").Value = "  +6.04%  "

$ws.Range("E42").Value = "  +1.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.34%  "

$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").Value = "1.959.15"
$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("E47").Value = "  -2.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.58%  "
